$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Swap "Recorded By" ordering wherever System + the user both recorded
#    attendance for a session ("System, x" -> "x, System").
# ---------------------------------------------------------------------------
$ws.Cells.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System")

# ---------------------------------------------------------------------------
# Helper: assign a percentage-looking string to a cell while keeping it as
# literal text (Excel would otherwise silently convert "78.3%" into the
# number 0.783 with a percent format). After writing the text we restore
# the original cell formatting (pulled from an untouched donor cell that
# still has the pristine style) so the cell's style index in the saved
# file is unaffected.
# ---------------------------------------------------------------------------
$fmtDonor = $ws.Range("M21")

function Set-TextValue($addr, $text) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $fmtDonor.Copy()
    $rng.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------------
# 2. Class Statistics summary table (K/L columns)
# ---------------------------------------------------------------------------
$ws.Range("L6").Value = 249
$ws.Range("L7").Value = 27
Set-TextValue "L9" "78.3%"

# ---------------------------------------------------------------------------
# 3. Group Statistics table rows 21-26 (O, P, R, S columns)
# ---------------------------------------------------------------------------
$ws.Range("O21").Value = 20
$ws.Range("P21").Value = 3
Set-TextValue "R21" "74.1%"
Set-TextValue "S21" "78.7%"

$ws.Range("O22").Value = 20
$ws.Range("P22").Value = 3
Set-TextValue "R22" "74.1%"
Set-TextValue "S22" "77.3%"

$ws.Range("O23").Value = 20
$ws.Range("P23").Value = 3
Set-TextValue "R23" "74.1%"
Set-TextValue "S23" "77.9%"

$ws.Range("O24").Value = 19
$ws.Range("P24").Value = 4
Set-TextValue "R24" "70.4%"
Set-TextValue "S24" "72.0%"

$ws.Range("O25").Value = 20
$ws.Range("P25").Value = 3
Set-TextValue "R25" "74.1%"
Set-TextValue "S25" "71.5%"

$ws.Range("O26").Value = 20
$ws.Range("P26").Value = 3
Set-TextValue "R26" "74.1%"
Set-TextValue "S26" "64.7%"

# ---------------------------------------------------------------------------
# 4. Rows 180, 207, 234, 261, 288, 315: the 14/01/2026 session became
#    "Recorded" for groups B1D1, B1D2, B1E1, B1E2, B1F1, B1F2.
#    Re-use the look of an already-"Recorded" row in the same group (green
#    fill / style) and then fill in the real attendance data.
# ---------------------------------------------------------------------------
$recordedRows = @(180, 207, 234, 261, 288, 315)
$sourceRows   = @(179, 206, 233, 260, 287, 314)
$hvalues      = @("19/23", "22/30", "17/26", "21/28", "20/26", "21/29")

for ($i = 0; $i -lt $recordedRows.Length; $i++) {
    $target = $recordedRows[$i]
    $source = $sourceRows[$i]

    $src = $ws.Range("A" + $source + ":I" + $source)
    $dst = $ws.Range("A" + $target + ":I" + $target)
    $src.Copy()
    $dst.PasteSpecial(-4122)

    $ws.Range("G" + $target).Value = "dnasr281@gmail.com"
    $ws.Range("H" + $target).Value = $hvalues[$i]
    $ws.Range("I" + $target).Value = "Recorded"
}
